# Applies the "major change" edit described in the commit:
# adds two new result columns (AG/AO blocks) for the (1000,10000) and
# (2000,10000) datasets to the existing row 31-38 table, and appends a
# brand-new results table for the (50,5000) dataset in rows 41-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AG31").Value = "(1000, 10000) "
$ws.Range("AO31").Value = "(2000, 10000) "
$ws.Range("AG32").Value = "data"
$ws.Range("AH32").Value = "k"
$ws.Range("AI32").Value = "Residual"
$ws.Range("AJ32").Value = "Objective"
$ws.Range("AK32").Value = "Time"
$ws.Range("AO32").Value = "data"
$ws.Range("AP32").Value = "k"
$ws.Range("AQ32").Value = "Residual"
$ws.Range("AR32").Value = "Objective"
$ws.Range("AS32").Value = "Time"
$ws.Range("AG33").Value = [double]"4"
$ws.Range("AH33").Value = [double]"494"
$ws.Range("AI33").Value = [double]"9.903114219051315e-07"
$ws.Range("AJ33").Value = [double]"1.818989403545856e-12"
$ws.Range("AK33").Value = [double]"69.00023085344583"
$ws.Range("AL33").Value = "NPG1"
$ws.Range("AO33").Value = [double]"4"
$ws.Range("AP33").Value = [double]"169"
$ws.Range("AQ33").Value = [double]"9.746984975540058e-07"
$ws.Range("AR33").Value = [double]"2.546585164964199e-11"
$ws.Range("AS33").Value = [double]"27.94447601307184"
$ws.Range("AT33").Value = "NPG1"
$ws.Range("AH34").Value = [double]"445"
$ws.Range("AI34").Value = [double]"8.768577252238355e-07"
$ws.Range("AJ34").Value = [double]"0"
$ws.Range("AK34").Value = [double]"64.28357172198594"
$ws.Range("AL34").Value = "NPG2"
$ws.Range("AP34").Value = [double]"153"
$ws.Range("AQ34").Value = [double]"9.797821944636511e-07"
$ws.Range("AR34").Value = [double]"1.091393642127514e-11"
$ws.Range("AS34").Value = [double]"24.39155614189804"
$ws.Range("AT34").Value = "NPG2"
$ws.Range("AH35").Value = [double]"1156"
$ws.Range("AI35").Value = [double]"9.996581219359656e-07"
$ws.Range("AJ35").Value = [double]"1.2732925824821e-11"
$ws.Range("AK35").Value = [double]"157.0114721460268"
$ws.Range("AL35").Value = "AdPG"
$ws.Range("AP35").Value = [double]"374"
$ws.Range("AQ35").Value = [double]"9.070735907268453e-07"
$ws.Range("AR35").Value = [double]"5.456968210637569e-12"
$ws.Range("AS35").Value = [double]"60.43496393319219"
$ws.Range("AT35").Value = "AdPG"
$ws.Range("AH36").Value = [double]"1199"
$ws.Range("AI36").Value = [double]"9.959513476882385e-07"
$ws.Range("AJ36").Value = [double]"1.2732925824821e-11"
$ws.Range("AK36").Value = [double]"166.1632217271253"
$ws.Range("AL36").Value = "AdaPG[1.5, 0.75]"
$ws.Range("AP36").Value = [double]"375"
$ws.Range("AQ36").Value = [double]"9.96443692078795e-07"
$ws.Range("AR36").Value = [double]"9.094947017729282e-12"
$ws.Range("AS36").Value = [double]"71.75540275964886"
$ws.Range("AT36").Value = "AdaPG[1.5, 0.75]"
$ws.Range("AH37").Value = [double]"2000"
$ws.Range("AI37").Value = [double]"75.98504085327562"
$ws.Range("AJ37").Value = [double]"1.982698449864984e-10"
$ws.Range("AK37").Value = [double]"319.7198628960177"
$ws.Range("AL37").Value = "PG-LS[1.1, 0.5]"
$ws.Range("AP37").Value = [double]"1000"
$ws.Range("AQ37").Value = [double]"3.447065752526245"
$ws.Range("AR37").Value = [double]"2.000888343900442e-11"
$ws.Range("AS37").Value = [double]"216.9363348800689"
$ws.Range("AT37").Value = "PG-LS[1.1, 0.5]"
$ws.Range("AH38").Value = [double]"2000"
$ws.Range("AI38").Value = [double]"0.005630165025612393"
$ws.Range("AJ38").Value = [double]"2.455635694786906e-10"
$ws.Range("AK38").Value = [double]"409.0998869882897"
$ws.Range("AL38").Value = "PG-LS[1.2, 0.5]"
$ws.Range("AP38").Value = [double]"1000"
$ws.Range("AQ38").Value = [double]"37.26707100255064"
$ws.Range("AR38").Value = [double]"5.275069270282984e-11"
$ws.Range("AS38").Value = [double]"243.1232928987592"
$ws.Range("AT38").Value = "PG-LS[1.2, 0.5]"
$ws.Range("A41").Value = "(50, 5000) "
$ws.Range("A42").Value = "data"
$ws.Range("B42").Value = "k"
$ws.Range("C42").Value = "Residual"
$ws.Range("D42").Value = "Objective"
$ws.Range("E42").Value = "Time"
$ws.Range("A43").Value = [double]"5"
$ws.Range("B43").Value = [double]"25370"
$ws.Range("C43").Value = [double]"9.997310441868259e-07"
$ws.Range("D43").Value = [double]"0"
$ws.Range("E43").Value = [double]"49.94762175902724"
$ws.Range("F43").Value = "NPG1"
$ws.Range("B44").Value = [double]"32926"
$ws.Range("C44").Value = [double]"9.999727942538714e-07"
$ws.Range("D44").Value = [double]"1.818989403545856e-12"
$ws.Range("E44").Value = [double]"40.21087073627859"
$ws.Range("F44").Value = "NPG2"
$ws.Range("B45").Value = [double]"50000"
$ws.Range("C45").Value = [double]"1.500767949450377e-06"
$ws.Range("D45").Value = [double]"1.731677912175655e-09"
$ws.Range("E45").Value = [double]"50.72289928887039"
$ws.Range("F45").Value = "AdPG"
$ws.Range("B46").Value = [double]"50000"
$ws.Range("C46").Value = [double]"1.429748190058816e-06"
$ws.Range("D46").Value = [double]"3.378772817086428e-09"
$ws.Range("E46").Value = [double]"63.33620479982346"
$ws.Range("F46").Value = "AdaPG[1.5, 0.75]"
$ws.Range("B47").Value = [double]"50000"
$ws.Range("C47").Value = [double]"1.976956301760311e-05"
$ws.Range("D47").Value = [double]"1.215397787746042e-06"
$ws.Range("E47").Value = [double]"54.56627183035016"
$ws.Range("F47").Value = "PG-LS[1.1, 0.5]"
$ws.Range("B48").Value = [double]"50000"
$ws.Range("C48").Value = [double]"1.836204009806177e-05"
$ws.Range("D48").Value = [double]"1.04966511571547e-06"
$ws.Range("E48").Value = [double]"66.18509535398334"
$ws.Range("F48").Value = "PG-LS[1.2, 0.5]"

Write-Host "Edit applied."
